# Weekly update: a new price-survey record for "Cilantro" (Macroferia
# Regional de Talca) was collected, so it is inserted as the new row 23.
# Every existing record from the old row 23 onward shifts down by one row
# (old row 40 becomes row 41), matching the target diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 23:40 down to 24:41, opening a blank row 23.
$ws.Rows.Item(23).Insert()

# Populate the newly opened row 23 with the new survey record.
$ws.Cells.Item(23, 1).Value  = 5
$ws.Cells.Item(23, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(23, 3).Value  = "Maule"
$ws.Cells.Item(23, 4).Value  = 44777
$ws.Cells.Item(23, 5).Value  = 7
$ws.Cells.Item(23, 6).Value  = 100112040
$ws.Cells.Item(23, 7).Value  = "Cilantro"
$ws.Cells.Item(23, 8).Value  = "Sin especificar"
$ws.Cells.Item(23, 9).Value  = "Primera"
$ws.Cells.Item(23, 10).Value = 150
$ws.Cells.Item(23, 11).Value = 13000
$ws.Cells.Item(23, 12).Value = 13000
$ws.Cells.Item(23, 13).Value = 13000
$ws.Cells.Item(23, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(23, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(23, 16).Value = 361
$ws.Cells.Item(23, 17).Value = 36
$ws.Cells.Item(23, 18).Value = "Hortaliza"

# Keep the date cell formatted the same way as the rest of column D.
$ws.Cells.Item(23, 4).NumberFormat = $ws.Cells.Item(24, 4).NumberFormat
